# Nuevo formato 15 jun 2021
# Updates "Aprobados/Reprobados/Por_Apro/Por_Repro/Promedio/Blancos/Por_Blan"
# statistics for several groups across the three partial-exam sheets.

$wb = $excel.ActiveWorkbook

# Columns: D=Totales, E=Aprobados, F=Reprobados, G=Por_Apro, H=Por_Repro,
#          I=Promedio, J=Blancos, K=Por_Blan

# --- Sheet "1er Parcial" ---
$ws = $wb.Worksheets.Item("1er Parcial")

$ws.Cells.Item(18, 5).Value = 33      # E18 Aprobados
$ws.Cells.Item(18, 6).Value = 0       # F18 Reprobados
$ws.Cells.Item(18, 7).Value = 100     # G18 Por_Apro
$ws.Cells.Item(18, 8).Value = 0       # H18 Por_Repro
$ws.Cells.Item(18, 9).Value = 7.8     # I18 Promedio
$ws.Cells.Item(18, 10).Value = 0      # J18 Blancos
$ws.Cells.Item(18, 11).Value = 0      # K18 Por_Blan

$ws.Cells.Item(19, 5).Value = 38
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 100
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 7.3
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0

$ws.Cells.Item(20, 5).Value = 38
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 97.44
$ws.Cells.Item(20, 8).Value = 2.56
$ws.Cells.Item(20, 9).Value = 7.7
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0

$ws.Cells.Item(21, 5).Value = 34
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 100
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 7.9
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0

# --- Sheet "2o Parcial" ---
$ws = $wb.Worksheets.Item("2o Parcial")

$ws.Cells.Item(18, 5).Value = 33
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 100
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 7.6
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0

$ws.Cells.Item(19, 5).Value = 38
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 100
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 7.1
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0

$ws.Cells.Item(20, 5).Value = 38
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 97.44
$ws.Cells.Item(20, 8).Value = 2.56
$ws.Cells.Item(20, 9).Value = 7.4
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0

$ws.Cells.Item(21, 5).Value = 34
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 100
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 7.6
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0

# --- Sheet "3er Parcial" ---
$ws = $wb.Worksheets.Item("3er Parcial")

$ws.Cells.Item(12, 5).Value = 32
$ws.Cells.Item(12, 6).Value = 4
$ws.Cells.Item(12, 7).Value = 88.89
$ws.Cells.Item(12, 8).Value = 11.11
$ws.Cells.Item(12, 9).Value = 6.7

$ws.Cells.Item(18, 5).Value = 33
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 100
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 7.7
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0

$ws.Cells.Item(19, 5).Value = 38
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 100
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 7.1
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0

$ws.Cells.Item(20, 5).Value = 38
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 97.44
$ws.Cells.Item(20, 8).Value = 2.56
$ws.Cells.Item(20, 9).Value = 7.4
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0

$ws.Cells.Item(21, 5).Value = 34
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 100
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 7.6
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
